$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fitness values (column C) for rows 2-49 from 7310 to 7318
$ws.Range("C2:C49").Value = 7318

# Update Fitness values (column C) for rows 206-252 to 7310
$ws.Range("C206:C252").Value = 7310
